$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.848.88"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.829.33"
$ws.Range("E3").Value = "  -1.66%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.41"
$ws.Range("E5").Value = "  +0.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6903"
$ws.Range("E6").Value = "  -1.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9995"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07662"
$ws.Range("E8").Value = "  -2.31%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3044"
$ws.Range("E9").Value = "  -2.46%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.21"
$ws.Range("E10").Value = "  -3.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07784"
$ws.Range("E11").Value = "  -0.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "92.90"
$ws.Range("E12").Value = "  +0.62%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.831.86"
$ws.Range("E13").Value = "  -0.80%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.080"
$ws.Range("E14").Value = "  -1.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6800"
$ws.Range("E15").Value = "  -1.80%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.439"
$ws.Range("E16").Value = "  -1.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008249"
$ws.Range("E17").Value = "  -2.76%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "28.842.55"
$ws.Range("E18").Value = "  -1.32%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.37"
$ws.Range("E19").Value = "  -2.78%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.073.40"
$ws.Range("E20").Value = "  -0.56%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.65"
$ws.Range("E21").Value = "  -2.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.438"
$ws.Range("E23").Value = "  -1.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9992"
$ws.Range("E24").Value = "  -0.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1494"
$ws.Range("E25").Value = "  -2.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.76"
$ws.Range("E26").Value = "  -0.48%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.725"
$ws.Range("E27").Value = "  -2.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.19"
$ws.Range("E28").Value = "  -2.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.538"
$ws.Range("E29").Value = "  -2.59%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.221"
$ws.Range("E30").Value = "  -1.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.154"
$ws.Range("E31").Value = "  -2.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.192"
$ws.Range("E32").Value = "  -1.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05110"
$ws.Range("E33").Value = "  -2.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7742"
$ws.Range("E34").Value = "  +3.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.854"
$ws.Range("E35").Value = "  -1.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.137"
$ws.Range("E36").Value = "  -3.06%  "

$ws.Range("E37").Value = "  -0.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.699"
$ws.Range("E40").Value = "  -1.62%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9570"
$ws.Range("E41").Value = "  +6.22%  "

$ws.Range("E42").Value = "  +2.35%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "106.85"
$ws.Range("E43").Value = "  -4.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9991"
$ws.Range("E44").Value = "  -0.17%  "

$ws.Range("E45").Value = "  +1.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5161"
$ws.Range("E46").Value = "  -0.39%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.973.93"
$ws.Range("E47").Value = "  -1.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.72"
$ws.Range("E48").Value = "  -7.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.748"
$ws.Range("E49").Value = "  -1.66%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000118"
$ws.Range("E50").Value = "  -5.37%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.928"
$ws.Range("E51").Value = "  -0.97%  "

# Row 38/39 content swap: VeChain now listed before Maker, with updated values
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01856"
$ws.Range("E38").Value = "  -0.42%  "

$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.265.58"
$ws.Range("E39").Value = "  +1.44%  "
